# Apply weekly update: shift rows 321-417 down to 323-419 (inserting 2 new
# rows), and populate the newly freed rows 321-322 with the latest week's
# observations for "Albahaca" at "Mercado Mayorista Lo Valledor de Santiago".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 321; this pushes the existing rows
# 321:417 down to 323:419 (matching the OOXML diff), carrying over their
# values/styles automatically.
$ws.Rows("321:322").Insert()

# Row 321 (new data)
$ws.Range("A321").Value2 = 6
$ws.Range("B321").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C321").Value2 = "Metropolitana"
$ws.Range("D321").Value2 = 44627
$ws.Range("E321").Value2 = 13
$ws.Range("F321").Value2 = 100112052
$ws.Range("G321").Value2 = "Albahaca"
$ws.Range("H321").Value2 = "Sin especificar"
$ws.Range("I321").Value2 = "Primera"
$ws.Range("J321").Value2 = 350
$ws.Range("K321").Value2 = 3000
$ws.Range("L321").Value2 = 3500
$ws.Range("M321").Value2 = 3143
$ws.Range("N321").Value2 = "`$/docena de matas"
$ws.Range("O321").Value2 = "Región Metropolitana"
$ws.Range("P321").Value2 = 524
$ws.Range("Q321").Value2 = 6
$ws.Range("R321").Value2 = "Hortaliza"

# Row 322 (new data)
$ws.Range("A322").Value2 = 6
$ws.Range("B322").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C322").Value2 = "Metropolitana"
$ws.Range("D322").Value2 = 44627
$ws.Range("E322").Value2 = 13
$ws.Range("F322").Value2 = 100112052
$ws.Range("G322").Value2 = "Albahaca"
$ws.Range("H322").Value2 = "Sin especificar"
$ws.Range("I322").Value2 = "Segunda"
$ws.Range("J322").Value2 = 170
$ws.Range("K322").Value2 = 2000
$ws.Range("L322").Value2 = 2500
$ws.Range("M322").Value2 = 2176
$ws.Range("N322").Value2 = "`$/docena de matas"
$ws.Range("O322").Value2 = "Región Metropolitana"
$ws.Range("P322").Value2 = 363
$ws.Range("Q322").Value2 = 6
$ws.Range("R322").Value2 = "Hortaliza"
